# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 422
$ws.Range("F4").Value = 162
$ws.Range("G5").Value = 178
$ws.Range("F6").Value = 3879
$ws.Range("F8").Value = 2551
$ws.Range("F10").Value = 3143
$ws.Range("F11").Value = 533
$ws.Range("F12").Value = 2318
$ws.Range("F15").Value = 93
$ws.Range("F16").Value = 453
$ws.Range("F19").Value = 210
$ws.Range("F20").Value = 348
$ws.Range("F22").Value = 399
$ws.Range("F23").Value = 661
$ws.Range("F27").Value = 1304
$ws.Range("F28").Value = 132
$ws.Range("F29").Value = 153
$ws.Range("F32").Value = 54
$ws.Range("F33").Value = 4316
$ws.Range("F34").Value = 4086
$ws.Range("F36").Value = 79
$ws.Range("F37").Value = 10
$ws.Range("F38").Value = 1130
$ws.Range("F40").Value = 475
$ws.Range("F42").Value = 1314
$ws.Range("F43").Value = 174
$ws.Range("F44").Value = 130
$ws.Range("F45").Value = 103
$ws.Range("F46").Value = 41
$ws.Range("F47").Value = 62
$ws.Range("F48").Value = 61

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F12").Value = 6
$ws.Range("F15").Value = 210

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2304

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 422
$ws.Range("F7").Value = 162
$ws.Range("G8").Value = 178
$ws.Range("F9").Value = 3879
$ws.Range("F11").Value = 2551
$ws.Range("F13").Value = 3143
$ws.Range("F14").Value = 533
$ws.Range("F15").Value = 2318
$ws.Range("F17").Value = 93
$ws.Range("F20").Value = 348
$ws.Range("F22").Value = 399
$ws.Range("F23").Value = 661
$ws.Range("F26").Value = 1304
$ws.Range("F27").Value = 132
$ws.Range("F28").Value = 153
$ws.Range("F30").Value = 54
$ws.Range("F32").Value = 4316
$ws.Range("F34").Value = 10
$ws.Range("F38").Value = 475
$ws.Range("F39").Value = 6
$ws.Range("F43").Value = 1314
$ws.Range("F44").Value = 174
$ws.Range("F45").Value = 103
$ws.Range("F46").Value = 41
$ws.Range("F47").Value = 62
$ws.Range("F48").Value = 61
$ws.Range("F49").Value = 210
